# "weekly box office data" - add a new "Weekly Data" worksheet that rolls
# the daily Box Office numbers up into NYT-bestseller-list weeks.

$wb = $excel.ActiveWorkbook

# -- NYT Books sheet: selection only moved -----------------------------
$wsBooks = $wb.Worksheets.Item("NYT Books")
$wsBooks.Range("D26").Select()

# -- Box Office sheet: scroll position changed, selection unchanged ----
$wsBoxOffice = $wb.Worksheets.Item("Box Office")
$wsBoxOffice.Activate()
try {
    $excel.ActiveWindow.ScrollRow = 56
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}
$wsBoxOffice.Range("A99").Select()

# -- New "Weekly Data" sheet, placed after "Box Office" -----------------
$new = $wb.Worksheets.Add()
$new.Name = "Weekly Data"
$new.Move($null, $wsBoxOffice)
$ws = $wb.Worksheets.Item("Weekly Data")

# Headers - write B1 before A1 so the shared-string table gets the same
# ordering ("Week End" before "Week Start") as the source workbook.
$ws.Range("B1").Value = "Week End"
$ws.Range("A1").Value = "Week Start"
$ws.Range("C1").Value = "Rank"
$ws.Range("D1").Value = "Box Office Gross"

# Week-end dates (column B) and bestseller rank (column C) for each of
# the 20 weekly rows.
$weekEnd = @(42966,42973,42980,42987,42994,43001,43008,43015,43022,43029,43036,43043,43050,43057,43064,43071,43078,43085,43092,43099)
$rank    = @(6,5,4,1,1,1,1,2,2,2,3,3,3,3,3,4,4,6,5,5)

for ($i = 0; $i -lt $weekEnd.Count; $i++) {
    $r = $i + 2
    $ws.Range("B$r").Value = $weekEnd[$i]
    $ws.Range("A$r").Formula = "=B$r-6"
    $ws.Range("C$r").Value = $rank[$i]
}

# Column D: weekly Box Office gross, summed from the daily "Box Office"
# sheet for the 15 weeks that are fully covered by that sheet's data.
$sumRanges = @{
    5  = "B2:B3"
    6  = "B4:B10"
    7  = "B11:B17"
    8  = "B18:B24"
    9  = "B25:B31"
    10 = "B32:B38"
    11 = "B39:B45"
    12 = "B46:B52"
    13 = "B53:B59"
    14 = "B60:B66"
    15 = "B67:B73"
    16 = "B74:B80"
    17 = "B81:B87"
    18 = "B88:B94"
    19 = "B95:B99"
}
foreach ($r in 5..19) {
    $ws.Range("D$r").Formula = "=SUM('Box Office'!$($sumRanges[$r]))"
}

# Number formats: dates in A/B (reuses the workbook's existing date
# format), currency in D (reuses the existing "$#,##0" format).
$ws.Range("A2:B21").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("D5:D19").NumberFormat = "`"$`"#,##0_);[Red]\(`"$`"#,##0\)"

# Column widths (best-fit-like) for the date and currency columns.
$ws.Columns("B:B").ColumnWidth = 13.666666666666666
$ws.Columns("D:D").ColumnWidth = 13.998697916666666

$ws.Range("D20").Select()

# New sheet becomes the active / selected tab, as the last one touched.
$ws.Activate()
